$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the selection to I1 (also clears the scrolled-away topLeftCell state
# that was left over at FC40/FE54 from the previous session).
$ws.Range("I1").Select()

# Force a full recalculation of the workbook, including volatile formulas
# (the DD/DG/DL/DO/DQ columns use RAND()-based formulas marked ca="1"),
# producing fresh cached values.
$excel.CalculateFull()
